# Create addAttachment as Admin user and TestSuccessfulAdminAddAttach
# Adds a new "testSuccessfulAdminAddAttachment" header/value block (rows 27-28,
# columns A-H) to the TestCaseDataSets sheet, below the existing
# testSuccessfulSumSheetDownLoadFiref block (rows 22-23), following the same
# header(yellow)/value/blank-rows pattern used throughout the sheet, and
# extends a couple of trailing blank formatted rows (34-35) to column H.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# --- 1. Build the new header row (27) and value row (28) ------------------
# Use the existing row16/row17 (header-style s=11 / data-style s=12) pair as a
# formatting template so the new rows get byte-identical cell styles, then
# overwrite the copied values with the new test's data.
$ws.Range("A16:E17").Copy($ws.Range("A27"))
$ws.Range("A16:C17").Copy($ws.Range("F27"))

# Header row 27: testname / client / username / password / mp / filetitle / documentpath / filename
$ws.Range("A27").Value = "testSuccessfulAdminAddAttachment"
$ws.Range("B27").Value = "client"
$ws.Range("C27").Value = "username"
$ws.Range("D27").Value = "password"
$ws.Range("E27").Value = "mp"
$ws.Range("F27").Value = "filetitle"
$ws.Range("G27").Value = "documentpath"
$ws.Range("H27").Value = "filename"

# Value row 28
$ws.Range("A28").Value = "testSuccessfulAdminAddAttachment"
$ws.Range("B28").Value = "Richmond County"
$ws.Range("C28").Value = "Jason Lee"
$ws.Range("D28").Value = "rese7"
$ws.Range("E28").Value = "152-3-182-00-0"
$ws.Range("F28").Value = "addingAttachment"
$ws.Range("G28").Value = "C:\\testfolder\\"
$ws.Range("H28").Value = "Bap1.JPEG"

# --- 2. Blank "custom formatted" separator rows around the new block ------
# Rows 26, 29-32 span the full A:H block; row 33 spans B:H only (mirrors the
# narrowing pattern already used by the legacy trailing rows).
$ws.Range("A26:H26").NumberFormat = "@"
$ws.Range("A26:H26").ClearContents()

$ws.Range("A29:H32").NumberFormat = "@"
$ws.Range("A29:H32").ClearContents()

$ws.Range("B33:H33").NumberFormat = "@"
$ws.Range("B33:H33").ClearContents()

# --- 3. Two fully-blank, formatted rows (34-35) extended to column H ------
$ws.Range("B34:H34").NumberFormat = "@"
$ws.Range("B34:H34").ClearContents()
$ws.Range("B35:H35").NumberFormat = "@"
$ws.Range("B35:H35").ClearContents()

# --- 4. Selection / view state --------------------------------------------
$ws.Range("D28").Select()
